$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: duplicate row 2's current content (region + kml path) BEFORE we overwrite row 2 ---
# (use .Text so we copy the literal displayed string, not a COM wrapper object)
$ws.Range("A3").Value = $ws.Range("A2").Text
$ws.Range("B3").Value = $ws.Range("B2").Text
$ws.Range("B3").Style = "Normal"

# --- Row 2: point at the new shapefile AOI path, and the first new test path ---
$ws.Range("B2").Value = "\\spatialfiles.bcgov\work\srm\nel\Local\Geomatics\Workarea\csostad\GitHubAutoAST\gss_authorizations\autoast\aoi.shp"
$ws.Range("F2").Value = "T:\test1"

# --- Row 3: second new test path ---
$ws.Range("F3").Value = "T:\test2"

# --- Row 8: empty cell pre-formatted with the built-in Hyperlink style (no content yet) ---
$ws.Hyperlinks.Add($ws.Range("B8"), "http://example.com")
$ws.Hyperlinks.Delete()
$ws.Range("B8").ClearContents()

# --- Update selection to reflect where the user ended up (B9) ---
$ws.Range("B9").Select()
